# Add a new "qPCR Primers" worksheet with a small table of qPCR primer
# sequences, matching the author's "Add files via upload" commit.

$wb = $excel.ActiveWorkbook

# Sheets we borrow existing cell formatting from, so the new sheet reuses
# the same style entries already present in the workbook instead of
# minting near-duplicate ones.
$wsNeuro   = $wb.Worksheets.Item("NeuropeptidesNewAndPublished")
$wsInSitu  = $wb.Worksheets.Item("InSituMarkers")

# New sheet goes at the very end of the tab strip.
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "qPCR Primers"

# ---- Header row -----------------------------------------------------
$ws.Range("A1").Value = "XLOC"
$ws.Range("B1").Value = "Primer"
$ws.Range("C1").Value = "Sequence"

# Bold Arial 10pt header style (becomes a brand-new font/cellXf entry,
# just like in the source workbook).
$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Name = "Arial"
$headerRange.Font.Size = 10
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 0

# ---- Column A (XLOC identifiers), filled top-to-bottom --------------
$ws.Range("A2").Value = "XLOC _006558"
$ws.Range("A3").Value = "XLOC_006558"
$ws.Range("A4").Value = "XLOC_006729"
$ws.Range("A5").Value = "XLOC_006729"
$ws.Range("A6").Value = "XLOC_000601"
$ws.Range("A7").Value = "XLOC_000601"
$ws.Range("A8").Value = "XLOC_008048"
$ws.Range("A9").Value = "XLOC_008048"

# ---- Rows 8-9 (Primer / Sequence), filled before rows 2-7 ------------
$ws.Range("B8").Value = "XLOC_008048_F_qPCR"
$ws.Range("C8").Value = "TGAAAAGtCTGTCCACACAATGGAAACC"
$ws.Range("B9").Value = "XLOC_008048_R_qPCR"
$ws.Range("C9").Value = "GCTGGTGCACAGTATGATGGACATG"

# ---- Rows 2-7 (Primer / Sequence) ------------------------------------
$ws.Range("B2").Value = "XLOC_006558_F_qPCR"
$ws.Range("C2").Value = "GACCTGGGTACGCTTACCTCAATAAAGG"
$ws.Range("B3").Value = "XLOC_006558_R_qPCR"
$ws.Range("C3").Value = "tgcTCCAAAGTTGCTCTATCAGGATGGT"
$ws.Range("B4").Value = "XLOC_006729_F_qPCR"
$ws.Range("C4").Value = "GGCGATTGGAGCCCGTTTGGTG"
$ws.Range("B5").Value = "XLOC_006729_R_qPCR"
$ws.Range("C5").Value = "ACAAACTTGGACCGTAAAACTGGT"
$ws.Range("B6").Value = "XLOC_000601_F_qPCR"
$ws.Range("C6").Value = "GAGGACCAAAACACGAAGCAGAAGATC"
$ws.Range("B7").Value = "XLOC_000601_R_qPCR"
$ws.Range("C7").Value = "TGGCATCTCTCCAGACAGGTTGG"

# ---- Row formatting ---------------------------------------------------
# Rows 2-7 are a little taller (16pt) than the sheet default.
$ws.Range("A2:C7").RowHeight = 16

# Column B / C styling on rows 2-3 and 6-7 (Calibri 12, theme-coloured
# text) mirrors formatting already used on the InSituMarkers sheet, so
# copy it across instead of re-creating it. (PasteSpecial only honours
# the first area of a multi-area range, so paste one cell at a time.)
foreach ($cell in @("B2", "B3", "B6", "B7")) {
    $wsInSitu.Range("C44").Copy() | Out-Null
    $ws.Range($cell).PasteSpecial(-4122) | Out-Null
}

# Column C on those same rows uses the plain "Sequence" style (Arial
# 10pt, theme colour) that already exists on NeuropeptidesNewAndPublished.
foreach ($cell in @("C2", "C3", "C6", "C7")) {
    $wsNeuro.Range("B1").Copy() | Out-Null
    $ws.Range($cell).PasteSpecial(-4122) | Out-Null
}

# Column B on rows 4-5 uses a Calibri 12pt / explicit black style that
# already exists on InSituMarkers.
foreach ($cell in @("B4", "B5")) {
    $wsInSitu.Range("B36").Copy() | Out-Null
    $ws.Range($cell).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = $false

# ---- View state --------------------------------------------------------
# NeuropeptidesNewAndPublished is scrolled down a bit.
$wsNeuro.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 10

# The new sheet becomes the active tab, with B7 selected - this also
# naturally clears the tabSelected flag that used to sit on InSituMarkers.
$ws.Activate() | Out-Null
$ws.Range("B7").Select() | Out-Null
